# Add a "Details" worksheet after "Schedule 1" in the Annual Report template,
# matching the header row / column layout of the new sheet, and add a frozen
# header row (pane split at row 1) to both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "Details" sheet right after "Schedule 1" ---------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Details"

# --- Populate header row values in the order that reproduces the target -
# shared-string table ordering (Debit, Credit, Period, Account Description,
# Act #1..Act #5 are new strings; Year/Amount reuse existing ones).
$ws2.Range("I1").Value = "Debit"
$ws2.Range("J1").Value = "Credit"
$ws2.Range("H1").Value = "Period"
$ws2.Range("B1").Value = "Account Description"
$ws2.Range("C1").Value = "Act #1"
$ws2.Range("D1").Value = "Act #2"
$ws2.Range("E1").Value = "Act #3"
$ws2.Range("F1").Value = "Act #4"
$ws2.Range("G1").Value = "Act #5"
$ws2.Range("A1").Value = "Year"
$ws2.Range("K1").Value = "Amount"

# --- Copy the header formatting from "Schedule 1" (same "Heading 2" style
# used for A1:F1, and the currency-formatted "Heading 2" style used for G1)
$ws1.Range("A1:F1").Copy()
$ws2.Range("A1:F1").PasteSpecial(-4122)

$ws1.Range("A1:B1").Copy()
$ws2.Range("G1:H1").PasteSpecial(-4122)

$ws1.Range("G1").Copy()
$ws2.Range("I1:K1").PasteSpecial(-4122)

# --- Column widths for the new sheet -------------------------------------
$ws2.Columns("A").ColumnWidth = 10.33203125
$ws2.Columns("B").ColumnWidth = 32.5546875
$ws2.Columns("C").ColumnWidth = 14.5546875
$ws2.Columns("D").ColumnWidth = 15.109375
$ws2.Columns("E").ColumnWidth = 14
$ws2.Columns("F").ColumnWidth = 14.109375
$ws2.Columns("G").ColumnWidth = 13.6640625
$ws2.Columns("H").ColumnWidth = 10.6640625
$ws2.Columns("I").ColumnWidth = 16.88671875
$ws2.Columns("J").ColumnWidth = 15.109375
$ws2.Columns("K").ColumnWidth = 20.109375

# --- Row 2 placeholder (keeps the thick-border bottom of the header row
# visually "open" under it, mirroring "Schedule 1"'s row 2) --------------
$ws2.Rows(2).RowHeight = 15

# --- Freeze the header row on both sheets --------------------------------
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Widen column G on "Schedule 1" and shrink the trailing filler-column
# block from H:U (8-21) down to H:T (8-20) -------------------------------
$ws1.Columns("G").ColumnWidth = 20.77734375
$ws1.Columns("U").ColumnWidth = 15.109375

# --- Re-select the original sheet / cell so the workbook opens the same -
$ws1.Activate()
$ws1.Range("A1").Select()
